$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated exchange-rate table: 매매기준율(B), 전일대비(C), 등락률(D)
# for rows 2-9 (미국 USD .. 뉴질랜드 NZD)

# B column: plain numeric 매매기준율 values
$ws.Cells.Item(2, 2).Value = 1188.3
$ws.Cells.Item(3, 2).Value = 1040.77
$ws.Cells.Item(4, 2).Value = 1340.82
$ws.Cells.Item(5, 2).Value = 185.87
$ws.Cells.Item(6, 2).Value = 1597.31
$ws.Cells.Item(7, 2).Value = 859.5599999999999
$ws.Cells.Item(8, 2).Value = 939.66
$ws.Cells.Item(9, 2).Value = 831.04

# C column: 전일대비 text values. The leading triangle glyph keeps Excel
# from treating these as numbers, so a direct assignment is safe.
$ws.Cells.Item(2, 3).Value = "▼1.70"
$ws.Cells.Item(3, 3).Value = "▼3.78"
$ws.Cells.Item(4, 3).Value = "▼4.06"
$ws.Cells.Item(5, 3).Value = "▼0.32"
$ws.Cells.Item(6, 3).Value = "▼1.34"
$ws.Cells.Item(7, 3).Value = "▼2.42"
$ws.Cells.Item(8, 3).Value = "▼1.31"
$ws.Cells.Item(9, 3).Value = "▼2.67"

# D column: 등락률 values that look like percentages (e.g. "-0.14%").
# A plain .Value assignment makes Excel auto-convert these into numeric
# percentages (with a new number format). To keep them as plain text -
# matching the source workbook, which stores them as shared-string text
# with the default/unstyled cell format - write them as a literal-string
# formula and then collapse the formula down to its computed value via
# copy/paste-values. That yields a plain text cell with no left-over
# number-format/style changes.
$dValues = @{
    2 = "-0.14%"
    3 = "-0.36%"
    4 = "-0.30%"
    5 = "-0.17%"
    6 = "-0.08%"
    7 = "-0.28%"
    8 = "-0.14%"
    9 = "-0.32%"
}

foreach ($r in ($dValues.Keys | Sort-Object)) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Formula = '="' + $dValues[$r] + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$excel.CutCopyMode = 0
